{"js": "// Load all paragraphs so we can locate the ones we need to touch.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) Insert a new \"Meta description\" paragraph right after the title\n//    (Heading1) paragraph at the top of the document.\n// ---------------------------------------------------------------------\nconst titlePara = paragraphs.items[0];\nconst insertionPoint = titlePara.getRange(Word.RangeLocation.after);\n\nconst metaDescriptionOoxml =\n  \"<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>\" +\n  \"<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>\" +\n  \"<pkg:xmlData>\" +\n  \"<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\" +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  \"<w:r/>\" +\n  \"<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>\" +\n  \"<w:r><w:t xml:space='preserve'>: Take a thrilling adventure to the underworld in this online slot game. Play free Book of Ba now and benefit from special expanding symbols and free spins.</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ninsertionPoint.insertOoxml(metaDescriptionOoxml, Word.InsertLocation.after);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Remove the duplicate \"Play Book of Ba Free...\" paragraph that used\n//    to sit near the end of the document (right before the closing\n//    italic paragraph), and update the italic paragraph's text with the\n//    new image-prompt copy.\n// ---------------------------------------------------------------------\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst duplicateTitleText =\n  \"Play Book of Ba Free - Exciting Egyptian-themed slot game\";\nconst oldClosingText =\n  \"Take a thrilling adventure to the underworld in this online slot game. Play free Book of Ba now and benefit from special expanding symbols and free spins.\";\n\nlet duplicateTitlePara = null;\nlet closingPara = null;\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const para = paragraphs.items[i];\n  if (closingPara === null && para.text === oldClosingText) {\n    closingPara = para;\n    continue;\n  }\n  if (duplicateTitlePara === null && para.text === duplicateTitleText) {\n    duplicateTitlePara = para;\n    break;\n  }\n}\n\nif (duplicateTitlePara) {\n  duplicateTitlePara.delete();\n  await context.sync();\n}\n\nconst newClosingText =\n  'Please create a feature image fitting the game \"Book of Ba\" with the following prompt: Create a cartoon-style image featuring a happy Maya warrior with glasses. The warrior should be holding a book in one hand and a torch in the other, standing in front of an entrance that leads to the underworld. Behind the warrior, there should be hieroglyphics and a glimpse of the riches that await in the darkness. The image should have a vibrant and exciting color scheme to match the thrill of the game.';\n\nif (closingPara) {\n  const found = closingPara.search(oldClosingText, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(newClosingText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Insert a new \"Meta description\" paragraph right after the title\n#    (Heading1) paragraph at the top of the document.\n# ---------------------------------------------------------------------\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n# Give the new paragraph body (Normal) formatting instead of inheriting\n# the Heading1 style of the paragraph it was split from.\n$metaPara.Range.Style = \"Normal\"\n\n$boldText = \"Meta description\"\n$restText = \": Take a thrilling adventure to the underworld in this online slot game. Play free Book of Ba now and benefit from special expanding symbols and free spins.\"\n\n$metaRange = $metaPara.Range\n# Exclude the trailing paragraph mark from the text we are about to set.\n$metaRange.MoveEnd(1, -1)\n$metaRange.Text = $boldText + $restText\n\n# Bold just the \"Meta description\" label.\n$boldRange = $d.Range($metaRange.Start, $metaRange.Start + $boldText.Length)\n$boldRange.Bold = 1\n\n# ---------------------------------------------------------------------\n# 2) Remove the duplicate \"Play Book of Ba Free...\" paragraph that used\n#    to sit near the end of the document (right before the closing\n#    italic paragraph), and update the italic paragraph's text with the\n#    new image-prompt copy.\n# ---------------------------------------------------------------------\n$duplicateTitleText = \"Play Book of Ba Free - Exciting Egyptian-themed slot game\"\n$oldClosingText = \"Take a thrilling adventure to the underworld in this online slot game. Play free Book of Ba now and benefit from special expanding symbols and free spins.\"\n$newClosingText = 'Please create a feature image fitting the game \"Book of Ba\" with the following prompt: Create a cartoon-style image featuring a happy Maya warrior with glasses. The warrior should be holding a book in one hand and a torch in the other, standing in front of an entrance that leads to the underworld. Behind the warrior, there should be hieroglyphics and a glimpse of the riches that await in the darkness. The image should have a vibrant and exciting color scheme to match the thrill of the game.'\n\n$dupIndex = -1\n$closingIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $paraText = $d.Paragraphs.Item($i).Range.Text\n  if ($paraText -eq ($duplicateTitleText + \"`r\") -and $i -ne 1) {\n    $dupIndex = $i\n  }\n  if ($paraText -eq ($oldClosingText + \"`r\")) {\n    $closingIndex = $i\n  }\n}\n\nif ($dupIndex -ne -1) {\n  $d.Paragraphs.Item($dupIndex).Range.Delete()\n  if ($closingIndex -gt $dupIndex) {\n    $closingIndex = $closingIndex - 1\n  }\n}\n\nif ($closingIndex -ne -1) {\n  $closingPara = $d.Paragraphs.Item($closingIndex)\n  $fullText = $closingPara.Range.Text\n  $startOffset = $fullText.IndexOf($oldClosingText)\n  if ($startOffset -ge 0) {\n    $rngStart = $closingPara.Range.Start + $startOffset\n    $rngEnd = $rngStart + $oldClosingText.Length\n    $targetRange = $d.Range($rngStart, $rngEnd)\n    $targetRange.Text = $newClosingText\n  }\n}\n"}
